$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.539699918141963
$ws.Range("D2").Value = 0.1554017772857748
$ws.Range("E2").Value = 0.1985850294174796
$ws.Range("F2").Value = 2.121403027924472
$ws.Range("G2").Value = 1.578716346805095
$ws.Range("H2").Value = 1.360084852487887
$ws.Range("J2").Value = 0.2803627659302492
$ws.Range("L2").Value = 0.45550550548964
$ws.Range("M2").Value = 0.423002409223642
$ws.Range("N2").Value = 1.949654226969095
# Row 3
$ws.Range("B3").Value = 1.48713522131581
$ws.Range("D3").Value = 0.149902513120125
$ws.Range("E3").Value = 0.1915218588740046
$ws.Range("F3").Value = 2.110718933671393
$ws.Range("G3").Value = 1.549299053078215
$ws.Range("H3").Value = 1.353321159804608
$ws.Range("J3").Value = 0.2703275671336769
$ws.Range("L3").Value = 0.422093885276638
$ws.Range("M3").Value = 0.4022333058380099
$ws.Range("N3").Value = 1.967734668824122
# Row 4
$ws.Range("B4").Value = 1.455500899976585
$ws.Range("D4").Value = 0.146483881268324
$ws.Range("E4").Value = 0.1871441150994713
$ws.Range("F4").Value = 2.10552084901019
$ws.Range("G4").Value = 1.532415551298783
$ws.Range("H4").Value = 1.349990450260776
$ws.Range("J4").Value = 0.264118042959808
$ws.Range("L4").Value = 0.4017128243601178
$ws.Range("M4").Value = 0.3896419721494269
$ws.Range("N4").Value = 1.979575509864148
# Row 5
$ws.Range("B5").Value = 1.44277129092896
$ws.Range("D5").Value = 0.1450800820317895
$ws.Range("E5").Value = 0.1853498347068765
$ws.Range("F5").Value = 2.10374419569213
$ws.Range("G5").Value = 1.525830467010252
$ws.Range("H5").Value = 1.348839454501615
$ws.Range("J5").Value = 0.2615756107209677
$ws.Range("L5").Value = 0.3934411584033057
$ws.Range("M5").Value = 0.3845515258338708
$ws.Range("N5").Value = 1.984586861980496
# Row 6
$ws.Range("B6").Value = 1.440667327608168
$ws.Range("D6").Value = 0.144846334279201
$ws.Range("E6").Value = 0.1850512722272768
$ws.Range("F6").Value = 2.103469790653932
$ws.Range("G6").Value = 1.524754801673311
$ws.Range("H6").Value = 1.348660778952706
$ws.Range("J6").Value = 0.2611527180915232
$ws.Range("L6").Value = 0.3920696990213628
$ws.Range("M6").Value = 0.3837087190621631
$ws.Range("N6").Value = 1.985430239716074
# Row 7
$ws.Range("B7").Value = 1.455328568849097
$ws.Range("D7").Value = 0.1464649925122359
$ws.Range("E7").Value = 0.1871199585953143
$ws.Range("F7").Value = 2.105495506452726
$ws.Range("G7").Value = 1.532325549568782
$ws.Range("H7").Value = 1.349974092866915
$ws.Range("J7").Value = 0.2640838033669795
$ws.Range("L7").Value = 0.4016011328171203
$ws.Range("M7").Value = 0.3895731559314228
$ws.Range("N7").Value = 1.979642340896518
# Row 8
$ws.Range("B8").Value = 1.521442956425801
$ws.Range("D8").Value = 0.1535142977612765
$ws.Range("E8").Value = 0.196158127287724
$ws.Range("F8").Value = 2.117435931793324
$ws.Range("G8").Value = 1.568327769147402
$ws.Range("H8").Value = 1.357581794141225
$ws.Range("J8").Value = 0.2769125841848705
$ws.Range("L8").Value = 0.4439574749586654
$ws.Range("M8").Value = 0.4158078896661053
$ws.Range("N8").Value = 1.955735019251328
# Row 9
$ws.Range("B9").Value = 1.656160549815183
$ws.Range("D9").Value = 0.1670093895985474
$ws.Range("E9").Value = 0.2135590207483844
$ws.Range("F9").Value = 2.151701943716958
$ws.Range("G9").Value = 1.64834852365459
$ws.Range("H9").Value = 1.379048087966595
$ws.Range("J9").Value = 0.3016895046370394
$ws.Range("L9").Value = 0.5280787067924564
$ws.Range("M9").Value = 0.4685282006153884
$ws.Range("N9").Value = 1.914711020442454
# Row 10
$ws.Range("B10").Value = 1.758218527721453
$ws.Range("D10").Value = 0.1767318673748122
$ws.Range("E10").Value = 0.2261502810200042
$ws.Range("F10").Value = 2.183559806445075
$ws.Range("G10").Value = 1.712984417820792
$ws.Range("H10").Value = 1.398847867726772
$ws.Range("J10").Value = 0.3196624420772025
$ws.Range("L10").Value = 0.5905346681691412
$ws.Range("M10").Value = 0.5080388349408906
$ws.Range("N10").Value = 1.888130585925573
# Row 11
$ws.Range("B11").Value = 1.805315530138898
$ws.Range("D11").Value = 0.1811148715438122
$ws.Range("E11").Value = 0.231837300299091
$ws.Range("F11").Value = 2.199519009645016
$ws.Range("G11").Value = 1.743680701866367
$ws.Range("H11").Value = 1.408738381787742
$ws.Range("J11").Value = 0.3277890830900674
$ws.Range("L11").Value = 0.6190908157021227
$ws.Range("M11").Value = 0.5261823861851553
$ws.Range("N11").Value = 1.876809197865846
# Row 12
$ws.Range("B12").Value = 1.82324600820499
$ws.Range("D12").Value = 0.1827690045391108
$ws.Range("E12").Value = 0.233985004339182
$ws.Range("F12").Value = 2.20577441415081
$ws.Range("G12").Value = 1.755492212137284
$ws.Range("H12").Value = 1.412611321871651
$ws.Range("J12").Value = 0.3308593285044736
$ws.Range("L12").Value = 0.6299250795857176
$ws.Range("M12").Value = 0.5330772575998353
$ws.Range("N12").Value = 1.872632688352276
# Row 13
$ws.Range("B13").Value = 1.819380103551509
$ws.Range("D13").Value = 0.182413005430405
$ws.Range("E13").Value = 0.233522717530235
$ws.Range("F13").Value = 2.204417755028018
$ws.Range("G13").Value = 1.752940028403827
$ws.Range("H13").Value = 1.411771529927904
$ws.Range("J13").Value = 0.3301984143091232
$ws.Range("L13").Value = 0.6275908096048965
$ws.Range("M13").Value = 0.5315912449867142
$ws.Range("N13").Value = 1.873527254470503
# Row 14
$ws.Range("B14").Value = 1.806788763050065
$ws.Range("D14").Value = 0.1812510699844836
$ws.Range("E14").Value = 0.2320141102537363
$ws.Range("F14").Value = 2.2000293906016
$ws.Range("G14").Value = 1.744648675346042
$ws.Range("H14").Value = 1.409054449959513
$ws.Range("J14").Value = 0.3280418173047934
$ws.Range("L14").Value = 0.6199817435376644
$ws.Range("M14").Value = 0.526749145130978
$ws.Range("N14").Value = 1.876463376761293
# Row 15
$ws.Range("B15").Value = 1.799088672507366
$ws.Range("D15").Value = 0.1805386228761421
$ws.Range("E15").Value = 0.2310892839733469
$ws.Range("F15").Value = 2.197369031156299
$ws.Range("G15").Value = 1.739594445632946
$ws.Range("H15").Value = 1.407406796641368
$ws.Range("J15").Value = 0.3267199086656802
$ws.Range("L15").Value = 0.6153236591742655
$ws.Range("M15").Value = 0.5237863799915772
$ws.Range("N15").Value = 1.878276244180078
# Row 16
$ws.Range("B16").Value = 1.755154071199229
$ws.Range("D16").Value = 0.1764446360434846
$ws.Range("E16").Value = 0.2257778015791203
$ws.Range("F16").Value = 2.182546443273765
$ws.Range("G16").Value = 1.711004495811494
$ws.Range("H16").Value = 1.398219326212399
$ws.Range("J16").Value = 0.3191303505089991
$ws.Range("L16").Value = 0.5886713559899874
$ws.Range("M16").Value = 0.5068565194100287
$ws.Range("N16").Value = 1.888885954839949
# Row 17
$ws.Range("B17").Value = 1.728372923525569
$ws.Range("D17").Value = 0.1739229902451882
$ws.Range("E17").Value = 0.2225089421414879
$ws.Range("F17").Value = 2.17382972585682
$ws.Range("G17").Value = 1.693797776451845
$ws.Range("H17").Value = 1.392809798893779
$ws.Range("J17").Value = 0.3144617351048282
$ws.Range("L17").Value = 0.5723579801100414
$ws.Range("M17").Value = 0.4965140307772131
$ws.Range("N17").Value = 1.895591859547579
# Row 18
$ws.Range("B18").Value = 1.713032239858876
$ws.Range("D18").Value = 0.1724688536309884
$ws.Range("E18").Value = 0.2206249382457131
$ws.Range("F18").Value = 2.168954108022788
$ws.Range("G18").Value = 1.684022544976983
$ws.Range("H18").Value = 1.389781505909838
$ws.Range("J18").Value = 0.3117718294911498
$ws.Range("L18").Value = 0.5629885769428995
$ws.Range("M18").Value = 0.4905813108654158
$ws.Range("N18").Value = 1.899521417299276
# Row 19
$ws.Range("B19").Value = 1.707849005752962
$ws.Range("D19").Value = 0.1719758601809929
$ws.Range("E19").Value = 0.2199863863270437
$ws.Range("F19").Value = 2.167326978354595
$ws.Range("G19").Value = 1.680733660689867
$ws.Range("H19").Value = 1.388770437286354
$ws.Range("J19").Value = 0.310860277799236
$ws.Range("L19").Value = 0.5598186010458335
$ws.Range("M19").Value = 0.4885753471697569
$ws.Range("N19").Value = 1.900864352856509
# Row 20
$ws.Range("B20").Value = 1.731217293161421
$ws.Range("D20").Value = 0.1741918113483649
$ws.Range("E20").Value = 0.2228573153263369
$ws.Range("F20").Value = 2.174743344446682
$ws.Range("G20").Value = 1.695616862913482
$ws.Range("H20").Value = 1.393377045208069
$ws.Range("J20").Value = 0.3149591978687454
$ws.Range("L20").Value = 0.5740931582163569
$ws.Range("M20").Value = 0.4976133508069083
$ws.Range("N20").Value = 1.894870501957257
# Row 21
$ws.Range("B21").Value = 1.810484547351052
$ws.Range("D21").Value = 0.1815925101314519
$ws.Range("E21").Value = 0.2324573832829699
$ws.Range("F21").Value = 2.201312597444058
$ws.Range("G21").Value = 1.74707894577918
$ws.Range("H21").Value = 1.409849054120372
$ws.Range("J21").Value = 0.3286754561535759
$ws.Range("L21").Value = 0.6222161506609041
$ws.Range("M21").Value = 0.5281707287593491
$ws.Range("N21").Value = 1.875597964073137
# Row 22
$ws.Range("B22").Value = 1.862848820454815
$ws.Range("D22").Value = 0.1863966136203885
$ws.Range("E22").Value = 0.2386975140963017
$ws.Range("F22").Value = 2.219913344274417
$ws.Range("G22").Value = 1.781805978506952
$ws.Range("H22").Value = 1.421358585475275
$ws.Range("J22").Value = 0.3375982042880992
$ws.Range("L22").Value = 0.6537877384636772
$ws.Range("M22").Value = 0.5482833327446315
$ws.Range("N22").Value = 1.863647186266356
# Row 23
$ws.Range("B23").Value = 1.834850086890356
$ws.Range("D23").Value = 0.1838355276269965
$ws.Range("E23").Value = 0.235370148765206
$ws.Range("F23").Value = 2.209872308533448
$ws.Range("G23").Value = 1.763170928535828
$ws.Range("H23").Value = 1.415147453718703
$ws.Range("J23").Value = 0.3328397875276181
$ws.Range("L23").Value = 0.6369264081967287
$ws.Range("M23").Value = 0.5375359466469689
$ws.Range("N23").Value = 1.869966559107169
# Row 24
$ws.Range("B24").Value = 1.72993117868856
$ws.Range("D24").Value = 0.1740702910712031
$ws.Range("E24").Value = 0.2226998304315231
$ws.Range("F24").Value = 2.174329874744998
$ws.Range("G24").Value = 1.69479408916277
$ws.Range("H24").Value = 1.393120338648032
$ws.Range("J24").Value = 0.3147343131946485
$ws.Range("L24").Value = 0.5733086549877555
$ws.Range("M24").Value = 0.4971163067370696
$ws.Range("N24").Value = 1.895196396451233
# Row 25
$ws.Range("B25").Value = 1.619174258298813
$ws.Range("D25").Value = 0.1633929687939144
$ws.Range("E25").Value = 0.2088858117215011
$ws.Range("F25").Value = 2.14126347280444
$ws.Range("G25").Value = 1.625681804781181
$ws.Range("H25").Value = 1.372536136007767
$ws.Range("J25").Value = 0.2950272870148041
$ws.Range("L25").Value = 0.5052075932231332
$ws.Range("M25").Value = 0.4541295426831056
$ws.Range("N25").Value = 1.925183153310876

Write-Host "Updated 240 cells in sheet pl_mw"
